# Applies the docxFiller data refresh described by the commit:
#   - replaces the liquidator's address block (two paragraphs) with the
#     new postal-code / street address wording, split into the same
#     run layout the live edit produced
#   - turns the hard-coded submission date into the "{5}" merge
#     placeholder, keeping the surrounding sentence intact
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph "пр-т Дмитра Яворницького, 37 " -> "49029" / ", м. Дніпро"
# (also drops the spell-check proofErr bookends that wrapped "пр-т")
# ---------------------------------------------------------------------
$addrPara = $d.Paragraphs.Item(9)

# Insert a clean paragraph ahead of it (inherits the same pPr/rPr) and
# delete the original paragraph outright so the stray <w:proofErr/>
# markers that flanked "пр-т" disappear with it.
$addrPara.Range.InsertParagraphBefore()
$oldAddrPara = $d.Paragraphs.Item(10)
$oldAddrPara.Range.Delete()

$newAddrPara = $d.Paragraphs.Item(9)
$rng = $newAddrPara.Range
$rng.MoveEnd(1, -1)
$rng.Text = "49029"

$tail = $rng.Duplicate
$tail.Collapse(0)
$tailStart = $tail.Start
$tail.InsertAfter(", м. Дніпро")

# Force the inserted tail into its own run (identical formatting, but a
# distinct <w:r>) by toggling a character property on and back off.
$splitRng = $d.Range($tailStart, $tailStart + 11)
$splitRng.Font.Bold = $true
$splitRng.Font.Bold = $false

# ---------------------------------------------------------------------
# Paragraph "49000, м. Дніпро" -> "вул." / " " / "Січових Стрільців, 28" / " "
# ---------------------------------------------------------------------
$streetPara = $d.Paragraphs.Item(10)
$rng2 = $streetPara.Range
$rng2.MoveEnd(1, -1)
$rng2.Text = "вул. Січових Стрільців, 28 "
$base2 = $rng2.Start

$seg1Len = 4    # "вул."
$seg2Len = 1    # " "
$seg3Len = 21   # "Січових Стрільців, 28"
$seg4Len = 1    # " "

$seg2 = $d.Range($base2 + $seg1Len, $base2 + $seg1Len + $seg2Len)
$seg2.Font.Bold = $true
$seg2.Font.Bold = $false

$seg3 = $d.Range($base2 + $seg1Len + $seg2Len, $base2 + $seg1Len + $seg2Len + $seg3Len)
$seg3.Font.Bold = $true
$seg3.Font.Bold = $false

$seg4 = $d.Range($base2 + $seg1Len + $seg2Len + $seg3Len, $base2 + $seg1Len + $seg2Len + $seg3Len + $seg4Len)
$seg4.Font.Bold = $true
$seg4.Font.Bold = $false

# ---------------------------------------------------------------------
# "... подано державному реєстратору 16 серпня 2018 року для ..."
#   -> "... подано державному реєстратору {5} року для ..."
# ---------------------------------------------------------------------
$dateRng = $d.Content
$dateRng.Find.Execute("16 серпня 2018") | Out-Null
$dateRng.Text = "{5}"
$dateRng.Font.Bold = $true
$dateRng.Font.Bold = $false

Write-Output "Applied docxFiller address/date updates"
